$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A77").Value = "agwatch_pdf_url"
$ws.Range("A78").Value = "agwatch_election_manifesto"
$ws.Range("B77").Value = "Manifesto URL on abgeordnetenwatch.de"
$ws.Range("B78").Value = "Is an electoral manifesto not just a general manifesto (AGWatch only)"

# Printer/page setup (matches paperSize/orientation emitted on save)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Scroll/selection state left by the editor after appending the new rows
$ws.Range("B81").Select() | Out-Null
